# Generate Report for Handoff
# The 30a2b24f-... entry moved from "Handed back" (already in sync) to
# "Ready for handoff" with updated handoff timestamps, and the three
# file rows were re-sorted (ffff5e5d, ffffffc6, 30a2b24f) in every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-25 10:34:14"

$ws1.Range("A3").Value = "ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-25 10:34:14"

$ws1.Range("A4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-25 10:38:06"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md", "", "", "ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md", "", "", "ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/30a2b24f-037d-49dc-813f-bef7f32643a6.md", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-25 10:34:04"
$ws2.Range("F2").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.md"
$ws2.Range("G2").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-25 10:34:44"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-25 10:34:04"
$ws2.Range("F3").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.md"
$ws2.Range("G3").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-25 10:34:44"
$ws2.Range("J3").Value = "Include"

$ws2.Range("A4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-25 10:37:57"
$ws2.Range("F4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.md"
$ws2.Range("G4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-03-25 10:36:59"
$ws2.Range("J4").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md", "", "", "ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2908635df447713fb8f3a865420e7d6f87dfa60f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/b817524bd32d0ca24a8c057fa4dbd3ee2f6bdf9f/e2e/7cd68582-18c7-4046-9456-dc3be6020c8c.md", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e17897311e3df4727b848c3be43d2c7cdca79abf/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md", "", "", "ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2908635df447713fb8f3a865420e7d6f87dfa60f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/b817524bd32d0ca24a8c057fa4dbd3ee2f6bdf9f/e2e/7cd68582-18c7-4046-9456-dc3be6020c8c.md", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e17897311e3df4727b848c3be43d2c7cdca79abf/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/30a2b24f-037d-49dc-813f-bef7f32643a6.md", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7aed6f6b117ce4d1ed4b0400594e2777d7456154/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.zh-cn.xlf", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/3d270b70bb1f8200a542d9d17646af566020be85/e2e/30a2b24f-037d-49dc-813f-bef7f32643a6.md", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.md")
$ws2.Hyperlinks.Add($ws2.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f380d4f23e4b1efb7b08638f5664aa8a4af6a434/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.zh-cn.xlf", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-25 10:34:14"
$ws3.Range("F2").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.md"
$ws3.Range("G2").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-25 10:34:59"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-25 10:34:14"
$ws3.Range("F3").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.md"
$ws3.Range("G3").Value = "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-25 10:34:59"
$ws3.Range("J3").Value = "Include"

$ws3.Range("A4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-25 10:38:06"
$ws3.Range("F4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.md"
$ws3.Range("G4").Value = "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.de-de.xlf"
$ws3.Range("H4").Value = "2016-03-25 10:37:14"
$ws3.Range("J4").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md", "", "", "ffff5e5d2b29-3177-41f6-bd2a-aa638438cf2d.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/129d5889f442f0cf26925f8ef61bf59c69cb9d98/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/1a2b1f2e1e9e616bdc2b3e5e37db5420214aa06b/e2e/7cd68582-18c7-4046-9456-dc3be6020c8c.md", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e9696f405de13aa1cf7d6e3b92fc61c23c9d3e85/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md", "", "", "ffffffc6bbfc5b-aea0-45dc-9294-c001807bed97.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/129d5889f442f0cf26925f8ef61bf59c69cb9d98/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/1a2b1f2e1e9e616bdc2b3e5e37db5420214aa06b/e2e/7cd68582-18c7-4046-9456-dc3be6020c8c.md", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e9696f405de13aa1cf7d6e3b92fc61c23c9d3e85/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf", "", "", "7cd68582-18c7-4046-9456-dc3be6020c8c.3e528f8f06461ac426c5f6a03c5c0fd62d6308ac.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/eea06f0ec2ec8538e1d872d1461e1da79d27c775/e2e/30a2b24f-037d-49dc-813f-bef7f32643a6.md", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b14efbd8116b615df730a9e7ba641236aa1caf7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.de-de.xlf", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/1b33f830307ea6a79a739f057e4dd816d23d3d4b/e2e/30a2b24f-037d-49dc-813f-bef7f32643a6.md", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.md")
$ws3.Hyperlinks.Add($ws3.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75bfbc7a015f95054f7508482e6b748ddc0e8da5/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.de-de.xlf", "", "", "30a2b24f-037d-49dc-813f-bef7f32643a6.28a88338cf71e2f8093030c2923088e5e4f697e7.de-de.xlf")

$wb.Save()
